$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("A22").Value = "TestCase_A21"
$ws.Range("C22").Value = "Verify View additional email preferences link is working"
$ws.Range("A23").Value = "TestCase_A22"
$ws.Range("A24").Value = "TestCase_A23"
$ws.Range("C24").Value = "Verify change password link in the account page is working correctly."
$ws.Range("C23").Value = 'Verify that the  checkbox  is present and can be modified for "Receive email notifications for likes,comments and other activity" is working correctly.'
$ws.Range("B23").Value = "OPQA-854,OPQA-853"
$ws.Range("B22").Value = "OPQA-399"
$ws.Range("B24").Value = "OPQA-527"

$ws.Range("D22").Value = "Y"
$ws.Range("E22").Value = "SKIP"
$ws.Range("D23").Value = "Y"
$ws.Range("E23").Value = "SKIP"
$ws.Range("D24").Value = "Y"
$ws.Range("E24").Value = "SKIP"

$ws.Range("A2").Copy()
$ws.Range("A22:E24").PasteSpecial(-4122)
$ws.Range("C22:C24").WrapText = $true
